$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the existing bordered header/data grid (A1:D5)
$ws.Range("A1:D5").HorizontalAlignment = -4108

# New EXTI pin-mapping column (entry order matters for shared-string table order)
$ws.Range("E1").Value = "EXTI8"
$ws.Range("E2").Value = "EXTI6"
$ws.Range("E5").Value = "EXTI7"
$ws.Range("E3").Value = "EXTI10"
$ws.Range("E4").Value = "EXTI9"
$ws.Range("E1:E5").HorizontalAlignment = -4108
$ws.Range("E1:E5").Borders.Color = 0
$ws.Range("E1:E5").Borders.LineStyle = 1

# New blank column F, centered
$ws.Range("F1:F6").HorizontalAlignment = -4108

# Row 6 (still unbordered) gets centered too
$ws.Range("A6:E6").HorizontalAlignment = -4108

$ws.Range("K17").Select()
